# Remove the two "faq" bullet paragraphs describing the question-button
# and its mail-answer behaviour, leaving the "Afisare lista de raspunsuri
# la intrebari frecvente" paragraph followed directly by the "Auth(..." one.

$d = $word.ActiveDocument

$wdParagraph = 4

# Locate the paragraph that starts with "-Buton care arata un formular..."
$startRange = $d.Content.Duplicate
$startRange.Find.Execute("-Buton care arata un formular", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startRange.Expand($wdParagraph)

# Locate the paragraph that contains "raspunsul la intrebare...prin mail"
$endRange = $d.Content.Duplicate
$endRange.Find.Execute("raspunsul la intrebare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endRange.Expand($wdParagraph)

# Build a range spanning both paragraphs (incl. their paragraph marks) and delete it
$delRange = $d.Range($startRange.Start, $endRange.End)
$delRange.Delete()
